# UseCaseDescriptions.docx edit:
# The merged "Extensions" row (gridSpan=2) in the "대여소 삭제" table is
# replaced by two normal (unmerged) rows:
#   Row 1: Cell1 = "2. 대여소 삭제"                         Cell2 = (empty)
#   Row 2: Cell1 = (empty)                                   Cell2 = "3. 대여소 삭제 완료 메시지 출력"

$d = $word.ActiveDocument

# --- locate the target table: the one whose 3rd row's first cell starts
#     with "Extensions" (6th table in the document) ---
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables($i)
    if ($candidate.Rows.Count -ge 3) {
        $probe = $candidate.Rows(3).Cells(1).Range.Text
        if ($probe -like "Extensions*") {
            $table = $candidate
        }
    }
}

$oldRowIndex = $table.Rows(3).Index

# Insert two fresh (unmerged, 2-cell) rows ahead of the old merged row.
# Add the "row that must end up second" first, then the "row that must end
# up first" -- Rows.Add(beforeRow) always inserts immediately above
# beforeRow, so inserting in this order yields the correct final order.
# NOTE: row/cell handles returned by this engine are index-based locators,
# not stable object identities -- always re-fetch the row by a freshly
# computed index rather than re-using a captured reference once the table
# shape has changed underneath it.
$rowDone = $table.Rows.Add($table.Rows($oldRowIndex))
$rowDone.Cells(2).Range.Text = "3. 대여소 삭제 완료 메시지 출력"

$rowNew = $table.Rows.Add($table.Rows($oldRowIndex))
$rowNew.Cells(1).Range.Text = "2. 대여소 삭제"

# Remove the original merged "Extensions" / "Step 1 ..." row entirely -- it
# has now been pushed two slots further down.
$table.Rows($oldRowIndex + 2).Delete()

Write-Output "Row split complete"
